$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 4, mirroring the structure of rows 2-3
$ws.Range("A4").Value = "test3"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "In Person"
$ws.Range("D4").Value = "All Over The World"
$ws.Range("E4").Value = "20"
$ws.Range("F4").Value = "15"
$ws.Range("G4").Value = "I'm currently busy"
$ws.Range("H4").Value = "10:08"
$ws.Range("I4").Value = "12:08"

# Update selection to match final state
$ws.Range("I5").Select()
